# Update countries & provincias Spain
# - Afganistan's case counts jumped, moving it up in the ranking to just
#   after Barein (row 60), pushing Kazajistan / Ghana / Moldavia down one
#   row each (their own figures are unchanged, only their position shifts).
# - Malasia (row 55) received refreshed figures.
# - The "datos actualizados" timestamp moved from 10:35 to 11:05.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Timestamp footer (A1)
$ws.Range("A1").Value = "Datos actualizados a 14 de Mayo de 2020 a las 11:05"

# Row 55 - Malasia: refreshed figures
$ws.Range("B55").Value = 6819
$ws.Range("C55").Value = 40
$ws.Range("D55").Value = 5351
$ws.Range("E55").Value = 1356
$ws.Range("F55").Value = 16
$ws.Range("G55").Value = 1
$ws.Range("H55").Value = 112

# Row 60 - now Afganistan (new, higher figures) - was Kazajistan
$ws.Range("A60").Value = "Afganistan"
$ws.Range("B60").Value = 5639
$ws.Range("C60").Value = 413
$ws.Range("D60").Value = 691
$ws.Range("E60").Value = 4812
$ws.Range("F60").Value = 7
$ws.Range("G60").Value = 4
$ws.Range("H60").Value = 136

# Row 61 - now Kazajistan (unchanged figures, shifted down one row) - was Ghana
$ws.Range("A61").Value = "Kazajistan"
$ws.Range("B61").Value = 5571
$ws.Range("C61").Value = 154
$ws.Range("D61").Value = 2408
$ws.Range("E61").Value = 3131
$ws.Range("F61").Value = 31
$ws.Range("G61").Value = 0
$ws.Range("H61").Value = 32

# Row 62 - now Ghana (unchanged figures, shifted down one row) - was Moldavia
$ws.Range("A62").Value = "Ghana"
$ws.Range("B62").Value = 5408
$ws.Range("C62").Value = 0
$ws.Range("D62").Value = 514
$ws.Range("E62").Value = 4870
$ws.Range("F62").Value = 5
$ws.Range("G62").Value = 0
$ws.Range("H62").Value = 24

# Row 63 - now Moldavia (unchanged figures, shifted down one row) - was Afganistan
$ws.Range("A63").Value = "Moldavia"
$ws.Range("B63").Value = 5406
$ws.Range("C63").Value = 0
$ws.Range("D63").Value = 2176
$ws.Range("E63").Value = 3045
$ws.Range("F63").Value = 251
$ws.Range("G63").Value = 0
$ws.Range("H63").Value = 185
